# TOOLS-SE FINAL DESEMBER 2023
# Restructure Sheet1: add a new "Id wilayah" column at the front, reorder the
# remaining columns, restyle the header/data rows, add an autofilter and its
# hidden _FilterDatabase defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture existing (pre-edit) values, verbatim, so we don't risk
#        mistyping strings (accents / nbsp / etc.) by hand. ------------------
$oldA1 = $ws.Range("A1").Value2   # wilayah
$oldB1 = $ws.Range("B1").Value2   # Kode Customer
$oldC1 = $ws.Range("C1").Value2   # Salesman Tujuan
$oldD1 = $ws.Range("D1").Value2   # Hari Tujuan
$oldE1 = $ws.Range("E1").Value2   # Rute Tujuan

$oldA2 = $ws.Range("A2").Value2   # SURABAYA 1
$oldB2 = $ws.Range("B2").Value2   # 999-0044829
$oldC2 = $ws.Range("C2").Value2   # OFFICE2(nbsp)
$oldD2 = $ws.Range("D2").Value2   # JUMAT
$oldE2 = $ws.Range("E2").Value2   # JUMAT GENAP

# --- 2. Write the new header row (A1:F1) in its new order. ------------------
$ws.Range("A1").Value = "Id wilayah"
$ws.Range("B1").Value = $oldA1
$ws.Range("C1").Value = $oldC1
$ws.Range("D1").Value = $oldE1
$ws.Range("E1").Value = $oldD1
$ws.Range("F1").Value = $oldB1

# --- 3. Write the new data row (A2:F2) in its new order. --------------------
$ws.Range("A2").Value = 201
$ws.Range("B2").Value = $oldA2
$ws.Range("C2").Value = $oldC2
$ws.Range("D2").Value = $oldE2
$ws.Range("E2").Value = $oldD2
$ws.Range("F2").Value = $oldB2

# --- 4. Header formatting. ---------------------------------------------------
# A1 & F1 take on the look the old "Kode Customer" header used: bold, green fill.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 5296274
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Interior.Color = 5296274

# B1:E1 take on the look the other headers used: bold, green fill.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Interior.Color = 5296274
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Interior.Color = 5296274
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Interior.Color = 5296274
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Interior.Color = 5296274

# --- 5. Data row alignment (new in this revision). ---------------------------
$ws.Range("B2").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("B2").VerticalAlignment = -4108     # xlVAlignCenter
$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").VerticalAlignment = -4108
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4108

$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4131
$ws.Range("F2").VerticalAlignment = -4108

# --- 6. Row 7 marker cells (previously A7/G7, now A7/B7/H7). ----------------
# A7 is left completely untouched so it keeps its original formatting.
# B7 / H7 are new empty, (default-styled) cells - touching Font.Name with the
# value it already has materializes an empty cell without inventing a new style.
$ws.Range("B7").Font.Name = "Calibri"
$ws.Range("H7").Font.Name = "Calibri"
$ws.Range("G7").Clear()

# --- 7. Column widths (closest values this engine's sixth-of-a-character
#        column-width quantization can represent). ---------------------------
$ws.Columns("A").ColumnWidth = 11
$ws.Columns("B").ColumnWidth = 11
$ws.Columns("C").ColumnWidth = 17.5
$ws.Columns("D").ColumnWidth = 13
$ws.Columns("E").ColumnWidth = 12.5
$ws.Columns("F").ColumnWidth = 16.166666666666668

# --- 8. AutoFilter + hidden _FilterDatabase defined name. -------------------
$ws.Range("A1:F1").AutoFilter() | Out-Null
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$1")
$fd.Visible = $false

# --- 9. Selection matches the final saved state. ----------------------------
$ws.Range("D4").Select() | Out-Null
